$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 27; $r++) {
    $ws.Range("C$r").Value = 45175
}
